# Revisão do Registar Stock
# - Rewords the two use-case steps in the "Cenário Normal" table
# - C7 now wraps onto multiple lines (style + row height) to fit the longer text
# - Selection/zoom brought in line with the saved view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reword the two scenario steps.
$ws.Range("C7").Value = "1. Categoriza e quantifica stock que chegou"
$ws.Range("D8").Value = "2. Adiciona stock ao sistema"

# 2. C7's new text is longer, so turn wrapping on and grow the row to fit it.
$ws.Range("C7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 38.25

# 3. Match the saved view state (zoom + active cell).
$excel.ActiveWindow.Zoom = 100
$ws.Range("C8").Select() | Out-Null
